$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluasi")

$ws.Range("B2").Value = 0.00516220261919316
$ws.Range("C2").Value = 1.374451949963991
$ws.Range("D2").Value = 0.767956424923307
$ws.Range("E2").Value = 0.00002664833588160472
$ws.Range("F2").Value = 30
$ws.Range("G2").Value = 1796.206321674443
$ws.Range("H2").Value = 42.38167436138458
$ws.Range("I2").Value = 24.32737954491605
$ws.Range("J2").Value = 0.7444787810482614
